# Auto-generated edit script: updates Hades_Profits price/profit columns (H-N)
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# the latest scheduled-runner price pull. Some rows gain/lose the optional
# LeveProfitNQ (M) / LeveProfitHQ (N) cell depending on whether an NQ/HQ price
# was available for that pull, matching upstream behavior.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 238.55556
$ws.Range("I53").Value = 100.666664
$ws.Range("J53").Value = 266.13333
$ws.Range("K53").Value = 100.666664
$ws.Range("L53").Value = 266.13333
$ws.Range("M53").Value = 536.333336
$ws.Range("N53").Value = -1540.13333
$ws.Range("H70").Value = 4621.2
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 4621.2
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 13863.6
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -14403.6
$ws.Range("H73").Value = 4621.2
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 4621.2
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 13863.6
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -15735.6
$ws.Range("H74").Value = 3975
$ws.Range("I74").Value = 3933.3333
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3933.3333
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -2997.3333
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 3975
$ws.Range("I77").Value = 3933.3333
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19666.6665
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -14986.6665
$ws.Range("N77").Value = -29360
$ws.Range("H129").Value = 846.5769
$ws.Range("J129").Value = 961.09753
$ws.Range("L129").Value = 2883.29259
$ws.Range("N129").Value = -12883.29259
$ws.Range("H137").Value = 2858879.2
$ws.Range("I137").Value = 4349283.5
$ws.Range("J137").Value = 2271.5833
$ws.Range("K137").Value = 13047850.5
$ws.Range("L137").Value = 6814.749899999999
$ws.Range("M137").Value = -13045300.5
$ws.Range("N137").Value = -11914.7499
$ws.Range("H138").Value = 2168679
$ws.Range("I138").Value = 2766.842
$ws.Range("J138").Value = 2878202
$ws.Range("K138").Value = 8300.526
$ws.Range("L138").Value = 8634606
$ws.Range("M138").Value = -3160.526
$ws.Range("N138").Value = -8644886

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1144.85
$ws.Range("I2").Value = 1061.9375
$ws.Range("J2").Value = 1476.5
$ws.Range("K2").Value = 1061.9375
$ws.Range("L2").Value = 1476.5
$ws.Range("M2").Value = -948.9375
$ws.Range("N2").Value = -1702.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H32").Value = 3804023
$ws.Range("I32").Value = 4129789.8
$ws.Range("J32").Value = 34436.43
$ws.Range("K32").Value = 4129789.8
$ws.Range("L32").Value = 34436.43
$ws.Range("M32").Value = -4129502.8
$ws.Range("N32").Value = -35010.43
$ws.Range("H116").Value = 1144.85
$ws.Range("I116").Value = 1061.9375
$ws.Range("J116").Value = 1476.5
$ws.Range("K116").Value = 1061.9375
$ws.Range("L116").Value = 1476.5
$ws.Range("M116").Value = 1232.0625
$ws.Range("N116").Value = -6064.5
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1189.3158
$ws.Range("I3").Value = 1112.7333
$ws.Range("J3").Value = 1476.5
$ws.Range("K3").Value = 1112.7333
$ws.Range("L3").Value = 1476.5
$ws.Range("M3").Value = -998.7333000000001
$ws.Range("N3").Value = -1704.5
$ws.Range("H94").Value = 1692
$ws.Range("I94").Value = 500
$ws.Range("K94").Value = 500
$ws.Range("M94").Value = -49

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7058.042
$ws.Range("I31").Value = 56691.2
$ws.Range("J31").Value = 1218.847
$ws.Range("K31").Value = 56691.2
$ws.Range("L31").Value = 1218.847
$ws.Range("M31").Value = -56396.2
$ws.Range("N31").Value = -1808.847
$ws.Range("H34").Value = 7058.042
$ws.Range("I34").Value = 56691.2
$ws.Range("J34").Value = 1218.847
$ws.Range("K34").Value = 56691.2
$ws.Range("L34").Value = 1218.847
$ws.Range("M34").Value = -56489.2
$ws.Range("N34").Value = -1622.847
$ws.Range("H62").Value = 9000
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 9666.666999999999
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 9666.666999999999
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -10914.667
$ws.Range("H65").Value = 9000
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 9666.666999999999
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 48333.335
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -54573.335
$ws.Range("H86").Value = 3310.5278
$ws.Range("I86").Value = 3163.4167
$ws.Range("J86").Value = 3604.75
$ws.Range("K86").Value = 3163.4167
$ws.Range("L86").Value = 3604.75
$ws.Range("M86").Value = -2040.4167
$ws.Range("N86").Value = -5850.75
$ws.Range("H89").Value = 3310.5278
$ws.Range("I89").Value = 3163.4167
$ws.Range("J89").Value = 3604.75
$ws.Range("K89").Value = 15817.0835
$ws.Range("L89").Value = 18023.75
$ws.Range("M89").Value = -10201.0835
$ws.Range("N89").Value = -29255.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H68").Value = 1165.1964
$ws.Range("I68").Value = 563.0323
$ws.Range("J68").Value = 1911.88
$ws.Range("K68").Value = 1689.0969
$ws.Range("L68").Value = 5735.64
$ws.Range("M68").Value = -878.0969
$ws.Range("N68").Value = -7357.64
$ws.Range("H71").Value = 1165.1964
$ws.Range("I71").Value = 563.0323
$ws.Range("J71").Value = 1911.88
$ws.Range("K71").Value = 5067.2907
$ws.Range("L71").Value = 17206.92
$ws.Range("M71").Value = -1011.2907
$ws.Range("N71").Value = -25318.92
$ws.Range("H107").Value = 861.89856
$ws.Range("J107").Value = 2229.0588
$ws.Range("L107").Value = 6687.176399999999
$ws.Range("N107").Value = -10527.1764
$ws.Range("H113").Value = 558.1875
$ws.Range("J113").Value = 577.7941
$ws.Range("L113").Value = 1733.3823
$ws.Range("N113").Value = -6073.382299999999
$ws.Range("H131").Value = 780.0540999999999
$ws.Range("J131").Value = 912.8077
$ws.Range("L131").Value = 2738.4231
$ws.Range("N131").Value = -12818.4231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9935.3125
$ws.Range("I80").Value = 27999.75
$ws.Range("J80").Value = 3913.8333
$ws.Range("K80").Value = 27999.75
$ws.Range("L80").Value = 3913.8333
$ws.Range("M80").Value = -27001.75
$ws.Range("N80").Value = -5909.8333
$ws.Range("H83").Value = 9935.3125
$ws.Range("I83").Value = 27999.75
$ws.Range("J83").Value = 3913.8333
$ws.Range("K83").Value = 139998.75
$ws.Range("L83").Value = 19569.1665
$ws.Range("M83").Value = -135006.75
$ws.Range("N83").Value = -29553.1665
$ws.Range("H102").Value = 1754.8928
$ws.Range("I102").Value = 1700.2632
$ws.Range("K102").Value = 1700.2632
$ws.Range("M102").Value = -78.2632000000001
$ws.Range("H122").Value = 3553
$ws.Range("I122").Value = 2850
$ws.Range("J122").Value = 4256
$ws.Range("K122").Value = 8550
$ws.Range("L122").Value = 12768
$ws.Range("M122").Value = -6100
$ws.Range("N122").Value = -17668
$ws.Range("H132").Value = 56797.46
$ws.Range("I132").Value = 44100.582
$ws.Range("J132").Value = 80237.84
$ws.Range("K132").Value = 132301.746
$ws.Range("L132").Value = 240713.52
$ws.Range("M132").Value = -129771.746
$ws.Range("N132").Value = -245773.52
$ws.Range("H135").Value = 58065.715
$ws.Range("J135").Value = 58065.715
$ws.Range("L135").Value = 58065.715
$ws.Range("N135").Value = -68205.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2319.077
$ws.Range("I93").Value = 2397.75
$ws.Range("J93").Value = 2193.2
$ws.Range("K93").Value = 2397.75
$ws.Range("L93").Value = 2193.2
$ws.Range("M93").Value = -1149.75
$ws.Range("N93").Value = -4689.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3467.3333
$ws.Range("I81").Value = 1501
$ws.Range("J81").Value = 3646.0908
$ws.Range("K81").Value = 3002
$ws.Range("L81").Value = 7292.1816
$ws.Range("M81").Value = -1941
$ws.Range("N81").Value = -9414.1816
$ws.Range("H84").Value = 3467.3333
$ws.Range("I84").Value = 1501
$ws.Range("J84").Value = 3646.0908
$ws.Range("K84").Value = 15010
$ws.Range("L84").Value = 36460.908
$ws.Range("M84").Value = -9706
$ws.Range("N84").Value = -47068.908
$ws.Range("H109").Value = 23547.2
$ws.Range("J109").Value = 23547.2
$ws.Range("L109").Value = 23547.2
$ws.Range("N109").Value = -26321.2

Write-Output "Updated 229 cells across 8 sheets"
